$d = $word.ActiveDocument

# --- Paragraph 1: title "Análise do Ciclo de Vida" ---
# Remove the two decorative rounded-rectangle background shapes that sit
# behind the title text (they were anchored drawings / AlternateContent).
$shapeCount = $d.Shapes.Count
for ($i = $shapeCount; $i -ge 1; $i--) {
    $shp = $d.Shapes.Item($i)
    if ($shp.Name -like "Retângulo*") {
        $shp.Delete()
    }
}

# Remove the centered alignment on the title paragraph (now left-aligned).
$p1 = $d.Paragraphs.Item(1)
$p1.Alignment = 0

# --- Paragraph 2: "Depósito de Dados Matrícula" -> split text, move bookmark ---
$p2 = $d.Paragraphs.Item(2)
$p2Start = $p2.Range.Start
$splitPoint = $p2Start + 17
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
